$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: quantity changed from 5 to 3 (total price recalculates automatically)
$ws.Range("B25").Value = 3

# Row 27: new component line item (N-Type low current MOSFET)
$ws.Range("A27").Value = 9845178
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = "N-Type low current MOSFET"
$ws.Range("D27").Value = 0.625
$ws.Range("F27").Value = "200mA, 60V"

# Row 48: header description wording fix
$ws.Range("C48").Value = "Single Row, 8 Pin male Headers"

# Row 49: header description wording fix
$ws.Range("C49").Value = "Single Row, 2 Pin male Header"

# Rows 55-57: new component line items in "other suppliers" table
$ws.Range("B55").Value = 1
$ws.Range("C55").Value = "100 ohm Resister"

$ws.Range("B56").Value = 1
$ws.Range("C56").Value = "30 pin, double row header"

$ws.Range("B57").Value = 1
$ws.Range("C57").Value = "42 pin, double row header"

# Update the active selection to reflect where the editor left off
$ws.Range("C57").Select()
